$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the cells whose values actually change are touched, matching the
# target diff exactly (row-wise permutation of the runs/balls/fours/sixes
# stats for Jonny Bairstow vs Sunrisers Hyderabad).

$ws.Range("C2").Value = "36"
$ws.Range("D2").Value = "28"
$ws.Range("E2").Value = "7"

$ws.Range("C3").Value = "97"
$ws.Range("D3").Value = "55"
$ws.Range("E3").Value = "7"
$ws.Range("F3").Value = "6"

$ws.Range("C4").Value = "53"
$ws.Range("D4").Value = "48"
$ws.Range("F4").Value = "1"

$ws.Range("C5").Value = "0"
$ws.Range("D5").Value = "3"
$ws.Range("E5").Value = "0"

$ws.Range("C6").Value = "10"
$ws.Range("D6").Value = "7"
$ws.Range("E6").Value = "1"

$ws.Range("C7").Value = "23"
$ws.Range("D7").Value = "24"
$ws.Range("E7").Value = "2"

$ws.Range("C8").Value = "16"
$ws.Range("D8").Value = "19"
$ws.Range("E8").Value = "0"
$ws.Range("F8").Value = "1"

$ws.Range("C9").Value = "19"
$ws.Range("D9").Value = "20"
$ws.Range("E9").Value = "4"
$ws.Range("F9").Value = "0"

$ws.Range("C10").Value = "61"
$ws.Range("D10").Value = "43"
$ws.Range("E10").Value = "6"
$ws.Range("F10").Value = "2"

$ws.Range("C11").Value = "25"
$ws.Range("D11").Value = "15"
$ws.Range("E11").Value = "2"
